$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before K (shifts K:AG -> L:AH), matching the
# "categories" column added to the Collection_AE specialization sheet.
$ws.Columns("K").Insert()
$ws.Range("K1").Value = "categories"

# The H4:H18 "Denormalized"/"Normalized" free-text cells lose their
# (no-op) explicit-font style, reverting to the default cell style.
$ws.Range("H4:H18").Style = "Normal"

# The column insert doesn't auto-grow the existing AutoFilter / hidden
# _FilterDatabase name, so extend both explicitly to the new last
# column (AH instead of AG).
$ws.AutoFilterMode = $false
$ws.Range("A1:AH34").AutoFilter()
$filterName = $wb.Names.Item(1)
$filterName.RefersTo = '=Collection_AE!$A$1:$AH$34'

# Match the author's final cursor/selection position.
$ws.Range("K20").Select()
